$d = $word.ActiveDocument

# Remove the empty paragraph, the "Ver no Jupiter..." paragraph and the
# "(c) 2020 ..." paragraph that follow the bibliography's closing line,
# while leaving the trailing blank paragraph (and the page-break
# paragraph after it) untouched. We match the paragraph marks with the
# wildcard token ^p so the whole paragraphs (including their marks) are
# removed and the remaining paragraphs reflow correctly.

$copyright = [char]169
$pattern = "^pVer no Jupiter Salvar em pdf Salvar em docx^p" + $copyright + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution^p"

$r = $d.Content.Duplicate
$r.Find.Execute($pattern, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
